$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRIDE_PROTEOMICS")

# Border: thin, color D2D2D2, applied to the full used range A1:L6
$full = $ws.Range("A1:L6")
$full.Borders.Color = 13816530
$full.Borders.LineStyle = 1

# Header-like style (row 1 + column A rows 2-6): bold white font, dark green fill
$headerRow = $ws.Range("A1:L1")
$headerRow.Font.Color = 16119285
$headerRow.Font.Bold = $true
$headerRow.Interior.Color = 4616993
$headerRow.VerticalAlignment = -4160

$headerCol = $ws.Range("A2:A6")
$headerCol.Font.Color = 16119285
$headerCol.Font.Bold = $true
$headerCol.Interior.Color = 4616993
$headerCol.VerticalAlignment = -4160

# Body style (B2:L6): light green fill
$bodyRng = $ws.Range("B2:L6")
$bodyRng.Interior.Color = 11783583
$bodyRng.VerticalAlignment = -4160
